$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value into a cell while forcing Excel to treat it as
# text (not auto-convert to a number), and then restore the cell's original
# (default/no-style) formatting so the saved style index doesn't change.
function Set-TextValue($cellRef, [string]$val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "68.498.31"
$ws.Range("E2").Value = "  -1.73%  "
Set-TextValue "D3" "2.444.68"
$ws.Range("E3").Value = "  -2.69%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue "D5" "562.49"
$ws.Range("E5").Value = "  -2.30%  "
Set-TextValue "D6" "162.25"
$ws.Range("E6").Value = "  -2.82%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("E9").Value = "  -6.91%  "
$ws.Range("E10").Value = "  -2.00%  "
Set-TextValue "D11" "0.341"
$ws.Range("E11").Value = "  -4.48%  "
Set-TextValue "D12" "4.80"
$ws.Range("E12").Value = "  -2.67%  "
Set-TextValue "D13" "2.896.19"
$ws.Range("E13").Value = "  -2.53%  "
Set-TextValue "D14" "68.388.34"
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("E15").Value = "  -4.32%  "
Set-TextValue "D16" "23.60"
$ws.Range("E16").Value = "  -5.33%  "
Set-TextValue "D17" "2.453.35"
$ws.Range("E17").Value = "  -1.83%  "
Set-TextValue "D18" "10.95"
$ws.Range("E18").Value = "  -3.13%  "
Set-TextValue "D19" "344.29"
$ws.Range("E19").Value = "  -1.81%  "
Set-TextValue "D20" "7.17"
$ws.Range("E20").Value = "  -4.94%  "
$ws.Range("E21").Value = "  -2.99%  "
Set-TextValue "D22" "0.999"
$ws.Range("E22").Value = "  -0.19%  "
Set-TextValue "D23" "1.87"
$ws.Range("E23").Value = "  -4.20%  "
Set-TextValue "D24" "68.17"
$ws.Range("E24").Value = "  -2.98%  "
Set-TextValue "D25" "3.75"
$ws.Range("E25").Value = "  -5.39%  "
$ws.Range("E26").Value = "  -2.04%  "
Set-TextValue "D27" "1.02"
$ws.Range("E27").Value = "  +2.28%  "
Set-TextValue "D28" "8.23"
$ws.Range("E28").Value = "  -7.13%  "
Set-TextValue "D29" "0.0₃0836"
$ws.Range("E29").Value = "  -6.51%  "
Set-TextValue "D30" "7.29"
$ws.Range("E30").Value = "  -7.22%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D31" "1.00"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D32" "1.17"
$ws.Range("E32").Value = "  -4.07%  "
Set-TextValue "D33" "429.12"
$ws.Range("E33").Value = "  -7.15%  "
Set-TextValue "D34" "1.67"
$ws.Range("E34").Value = "  -3.45%  "
$ws.Range("B35").Value = "POPCAT"
$ws.Range("C35").Value = "https://coinranking.com/coin/sLBuDEsp6+popcat-popcat"
Set-TextValue "D35" "2.98"
$ws.Range("E35").Value = "  +100.46%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D36" "156.42"
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").Value = "  -0.01%  "
Set-TextValue "D39" "0.109"
$ws.Range("E39").Value = "  -6.21%  "
Set-TextValue "D40" "17.89"
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("E41").Value = "  -4.30%  "
Set-TextValue "D42" "4.49"
$ws.Range("E42").Value = "  -4.38%  "
Set-TextValue "D43" "1.53"
$ws.Range("E43").Value = "  -4.70%  "
$ws.Range("E44").Value = "  +0.23%  "
Set-TextValue "D45" "2.08"
$ws.Range("E45").Value = "  -6.66%  "
Set-TextValue "D46" "134.85"
$ws.Range("E46").Value = "  -5.44%  "
Set-TextValue "D47" "3.36"
$ws.Range("E47").Value = "  -3.34%  "
Set-TextValue "D48" "0.490"
$ws.Range("E48").Value = "  -6.09%  "
$ws.Range("E49").Value = "  -2.19%  "
Set-TextValue "D50" "0.562"
$ws.Range("E50").Value = "  -3.34%  "
Set-TextValue "D51" "0.0914"
$ws.Range("E51").Value = "  -1.60%  "
